# [ADDITIONAL SCRAPING] added code to scrape more data about a player's
# batting performance in a match, also updated the excel sheets.
#
# 1. Turn the MATCH_CARD_LINK columns (ODI Batting / ODI Bowling) into
#    MATCH_CODE columns holding just the numeric match code instead of the
#    full howstat.com scorecard URL.
# 2. Add a new first sheet "Player Info" with the scraped player's basic
#    info (id / name / batting hand / bowling style).

$wb = $excel.ActiveWorkbook

function Set-TextValue($range, $text) {
    # Force the cell to be stored as text even when the text looks like a
    # plain number (e.g. "4259" or "4506") - mirrors typing `'4259` into
    # Excel - then drop back to the Normal style so no stray number format
    # sticks to the cell.
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

# ---------------------------------------------------------------------
# ODI Batting: MATCH_CARD_LINK -> MATCH_CODE
# ---------------------------------------------------------------------
$batting = $wb.Worksheets.Item("ODI Batting")
$battingRows = $batting.UsedRange.Rows.Count

$batting.Range("D1").Value = "MATCH_CODE"

for ($r = 2; $r -le $battingRows; $r++) {
    $cell = $batting.Cells.Item($r, 4)
    $link = $cell.Value2
    if ($link -and ($link -match "MatchCode=(\d+)")) {
        Set-TextValue $cell $matches[1]
    }
}

# ---------------------------------------------------------------------
# ODI Bowling: MATCH_CARD_LINK -> MATCH_CODE
# ---------------------------------------------------------------------
$bowling = $wb.Worksheets.Item("ODI Bowling")
$bowlingRows = $bowling.UsedRange.Rows.Count

$bowling.Range("B1").Value = "MATCH_CODE"

for ($r = 2; $r -le $bowlingRows; $r++) {
    $cell = $bowling.Cells.Item($r, 2)
    $link = $cell.Value2
    if ($link -and ($link -match "MatchCode=(\d+)")) {
        Set-TextValue $cell $matches[1]
    }
}

# ---------------------------------------------------------------------
# New "Player Info" sheet, inserted as the first sheet in the workbook.
# ---------------------------------------------------------------------
$info = $wb.Worksheets.Add()
$info.Name = "Player Info"

$info.Range("A1").Value = "ID"
$info.Range("B1").Value = "NAME"
$info.Range("C1").Value = "BATTING_HAND"
$info.Range("D1").Value = "BOWL_STYLE"

$header = $info.Range("A1:D1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160
$header.Borders.LineStyle = 1

Set-TextValue $info.Range("A2") "4506"
$info.Range("B2").Value = "Sayed Ahmad Shirzad"
$info.Range("C2").Value = "Left Handed"
$info.Range("D2").Value = "Left Arm Medium"

Write-Output "Player Info sheet added; MATCH_CARD_LINK columns converted to MATCH_CODE."
